# Generate Report for Handoff
# - Flip "Handed back: in sync with en-US" status to "Ready for handoff"
#   on every sheet that shows it, and refresh the associated handoff
#   timestamp cells.
# - The Status columns are narrower now that the text is shorter, so
#   shrink them to match.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status cells + "Latest HO Xliff Generate Date"
$ws_overview.Range("E2").Value = "Ready for handoff"
$ws_overview.Range("F2").Value = "Ready for handoff"
$ws_overview.Range("G2").Value = "2016-08-12 21:13:01"

# zh-cn detail sheet: Status + Latest Handoff Datetime
$ws_zhcn.Range("C2").Value = "Ready for handoff"
$ws_zhcn.Range("H2").Value = "2016-08-12 21:12:53"

# de-de detail sheet: Status + Latest Handoff Datetime
$ws_dede.Range("C2").Value = "Ready for handoff"
$ws_dede.Range("H2").Value = "2016-08-12 21:13:01"

# Re-fit the (now narrower) status columns on all three sheets to match
# the shorter "Ready for handoff" text.
$ws_overview.Range("E1:F1").EntireColumn.ColumnWidth = 16.33
$ws_zhcn.Range("C1").EntireColumn.ColumnWidth = 16.33
$ws_dede.Range("C1").EntireColumn.ColumnWidth = 16.33
